$wb = $excel.ActiveWorkbook

# ===== ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1632.2222
$ws.Range("I40").Value = 1199
$ws.Range("K40").Value = 1199
$ws.Range("M40").Value = -1024

# Row 80
$ws.Range("H80").Value = 811.9286
$ws.Range("I80").Value = 827.3125
$ws.Range("J80").Value = 791.4167
$ws.Range("K80").Value = 2481.9375
$ws.Range("L80").Value = 2374.2501
$ws.Range("M80").Value = -1483.9375
$ws.Range("N80").Value = -4370.2501

# Row 83
$ws.Range("H83").Value = 811.9286
$ws.Range("I83").Value = 827.3125
$ws.Range("J83").Value = 791.4167
$ws.Range("K83").Value = 7445.8125
$ws.Range("L83").Value = 7122.7503
$ws.Range("M83").Value = -2453.8125
$ws.Range("N83").Value = -17106.7503

# Row 86
$ws.Range("H86").Value = 2759.125
$ws.Range("I86").Value = 2848.25
$ws.Range("J86").Value = 2670
$ws.Range("K86").Value = 2848.25
$ws.Range("L86").Value = 2670
$ws.Range("M86").Value = -1725.25
$ws.Range("N86").Value = -4916

# Row 89
$ws.Range("H89").Value = 2759.125
$ws.Range("I89").Value = 2848.25
$ws.Range("J89").Value = 2670
$ws.Range("K89").Value = 14241.25
$ws.Range("L89").Value = 13350
$ws.Range("M89").Value = -8625.25
$ws.Range("N89").Value = -24582

# Row 132
$ws.Range("H132").Value = 7248997.5
$ws.Range("I132").Value = 9260736
$ws.Range("K132").Value = 27782208
$ws.Range("M132").Value = -27779678

# Row 138
$ws.Range("H138").Value = 6947435
$ws.Range("I138").Value = 983.3
$ws.Range("J138").Value = 11909186
$ws.Range("K138").Value = 2949.9
$ws.Range("L138").Value = 35727558
$ws.Range("M138").Value = 2190.1
$ws.Range("N138").Value = -35737838

# ===== ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 125008440
$ws.Range("I2").Value = 166672910
$ws.Range("K2").Value = 166672910
$ws.Range("M2").Value = -166672797

# Row 32
$ws.Range("H32").Value = 4755.672
$ws.Range("I32").Value = 2230.9482
$ws.Range("J32").Value = 29161.334
$ws.Range("K32").Value = 2230.9482
$ws.Range("L32").Value = 29161.334
$ws.Range("M32").Value = -1943.9482
$ws.Range("N32").Value = -29735.334

# Row 61
$ws.Range("H61").Value = 4986.951
$ws.Range("I61").Value = 4659.1763
$ws.Range("J61").Value = 5219.125
$ws.Range("K61").Value = 4659.1763
$ws.Range("L61").Value = 5219.125
$ws.Range("M61").Value = -4447.1763
$ws.Range("N61").Value = -5643.125

# Row 74
$ws.Range("H74").Value = 21068.316
$ws.Range("I74").Value = 20947.195
$ws.Range("J74").Value = 21754.666
$ws.Range("K74").Value = 20947.195
$ws.Range("L74").Value = 21754.666
$ws.Range("M74").Value = -20073.195
$ws.Range("N74").Value = -23502.666

# Row 77
$ws.Range("H77").Value = 21068.316
$ws.Range("I77").Value = 20947.195
$ws.Range("J77").Value = 21754.666
$ws.Range("K77").Value = 104735.975
$ws.Range("L77").Value = 108773.33
$ws.Range("M77").Value = -100367.975
$ws.Range("N77").Value = -117509.33

# Row 97
$ws.Range("H97").Value = 2014.75
$ws.Range("I97").Value = 1899.0769
$ws.Range("K97").Value = 1899.0769
$ws.Range("M97").Value = -1403.0769

# Row 110
$ws.Range("H110").Value = 6324
$ws.Range("I110").Value = 6066.476
$ws.Range("K110").Value = 6066.476
$ws.Range("M110").Value = -4021.476

# Row 116
$ws.Range("H116").Value = 125008440
$ws.Range("I116").Value = 166672910
$ws.Range("K116").Value = 166672910
$ws.Range("M116").Value = -166670616

# Row 136
$ws.Range("H136").Value = 4986.951
$ws.Range("I136").Value = 4659.1763
$ws.Range("J136").Value = 5219.125
$ws.Range("K136").Value = 13977.5289
$ws.Range("L136").Value = 15657.375
$ws.Range("M136").Value = -11427.5289
$ws.Range("N136").Value = -20757.375

# ===== BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 125008440
$ws.Range("I3").Value = 166672910
$ws.Range("K3").Value = 166672910
$ws.Range("M3").Value = -166672796

# Row 86
$ws.Range("H86").Value = 2570.6
$ws.Range("I86").Value = 2219.2727
$ws.Range("K86").Value = 2219.2727
$ws.Range("M86").Value = -1096.2727

# Row 89
$ws.Range("H89").Value = 2570.6
$ws.Range("I89").Value = 2219.2727
$ws.Range("K89").Value = 11096.3635
$ws.Range("M89").Value = -5480.363499999999

# Row 94
$ws.Range("H94").Value = 2308.05
$ws.Range("I94").Value = 2424.8
$ws.Range("J94").Value = 1957.8
$ws.Range("K94").Value = 2424.8
$ws.Range("L94").Value = 1957.8
$ws.Range("M94").Value = -1973.8
$ws.Range("N94").Value = -2859.8

# Row 107
$ws.Range("H107").Value = 9676.691999999999
$ws.Range("I107").Value = 1354.125
$ws.Range("K107").Value = 1354.125
$ws.Range("M107").Value = 565.875

# ===== CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 1829.5883
$ws.Range("J107").Value = 2285.7144
$ws.Range("L107").Value = 2285.7144
$ws.Range("N107").Value = -6125.7144

# Row 122
$ws.Range("H122").Value = 1543.3636
$ws.Range("I122").Value = 974.75
$ws.Range("J122").Value = 1868.2858
$ws.Range("K122").Value = 2924.25
$ws.Range("L122").Value = 5604.857400000001
$ws.Range("M122").Value = -474.25
$ws.Range("N122").Value = -10504.8574

# Row 134
$ws.Range("H134").Value = 21112.158
$ws.Range("I134").Value = 9010.143
$ws.Range("K134").Value = 27030.429
$ws.Range("M134").Value = -24495.429

# ===== CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 51
$ws.Range("H51").Value = 800.6
$ws.Range("I51").Value = 800.6
$ws.Range("K51").Value = 2401.8
$ws.Range("M51").Value = -1941.8

# ===== GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 32262660
$ws.Range("I102").Value = 1465.1724
$ws.Range("K102").Value = 1465.1724
$ws.Range("M102").Value = 156.8276000000001

# Row 113
$ws.Range("H113").Value = 1855.5555
$ws.Range("I113").Value = 1855.5555
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1855.5555
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 314.4445000000001
$ws.Range("N113").ClearContents()

# Row 126
$ws.Range("H126").Value = 17976.268
$ws.Range("I126").Value = 18903.143
$ws.Range("K126").Value = 56709.429
$ws.Range("M126").Value = -54239.429

# Row 132
$ws.Range("H132").Value = 2035.3182
$ws.Range("I132").Value = 2022.2407
$ws.Range("J132").Value = 2094.1667
$ws.Range("K132").Value = 6066.7221
$ws.Range("L132").Value = 6282.500100000001
$ws.Range("M132").Value = -3536.7221
$ws.Range("N132").Value = -11342.5001

# ===== LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4477.206
$ws.Range("I40").Value = 3816.6538
$ws.Range("J40").Value = 6624
$ws.Range("K40").Value = 3816.6538
$ws.Range("L40").Value = 6624
$ws.Range("M40").Value = -3680.6538
$ws.Range("N40").Value = -6896

# Row 132
$ws.Range("H132").Value = 2554.6943
$ws.Range("I132").Value = 2398.3125
$ws.Range("J132").Value = 3805.75
$ws.Range("K132").Value = 7194.9375
$ws.Range("L132").Value = 11417.25
$ws.Range("M132").Value = -4664.9375
$ws.Range("N132").Value = -16477.25

# ===== WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 94
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

# Row 96
$ws.Range("H96").Value = 1595.6154
$ws.Range("I96").Value = 1386.091
$ws.Range("K96").Value = 1386.091
$ws.Range("M96").Value = -13.09099999999989

# Row 97
$ws.Range("H97").Value = 50000
$ws.Range("J97").Value = 50000
$ws.Range("L97").Value = 50000
$ws.Range("N97").Value = -51982

# Row 107
$ws.Range("H107").Value = 219.44444
$ws.Range("J107").Value = 179.5
$ws.Range("L107").Value = 538.5
$ws.Range("N107").Value = -4378.5

# Row 113
$ws.Range("H113").Value = 1497.4286
$ws.Range("I113").Value = 1083.1428
$ws.Range("K113").Value = 3249.4284
$ws.Range("M113").Value = -1079.4284

# Row 122
$ws.Range("H122").Value = 1481.3529
$ws.Range("I122").Value = 1295.4
$ws.Range("J122").Value = 2876
$ws.Range("K122").Value = 3886.2
$ws.Range("L122").Value = 8628
$ws.Range("M122").Value = -1436.2
$ws.Range("N122").Value = -13528
